$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price (D) text looks numeric need the column forced to Text format first,
# otherwise Excel auto-converts the entered string to a float and mangles formatting
# (trailing zeros, float precision, etc.) exactly as it would for manual entry.
$textRows = @(5,10,15,16,19,20,25,35,43,44,45,47,48,50)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Per-cell Price (D) and Volume(1h) (E) updates
$ws.Range("D2").Value = "26.354.02"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").Value = "1.592.63"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "209.97"
$ws.Range("E5").Value = "  -0.76%  "
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("D10").Value = "19.55"
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "1.817.00"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "1.611.87"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "64.61"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "26.356.43"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("D19").Value = "7.50"
$ws.Range("E19").Value = "  +4.95%  "
$ws.Range("D20").Value = "211.60"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("E23").Value = "  -3.28%  "
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "145.36"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "1.303.33"
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("D35").Value = "0.614"
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  -13.19%  "
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("E42").Value = "  +3.09%  "
$ws.Range("D43").Value = "62.73"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "2.13"
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("D45").Value = "0.763"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").Value = "1.728.76"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "88.14"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D50").Value = "0.0984"
$ws.Range("E50").Value = "  -4.25%  "
$ws.Range("E51").Value = "  -1.41%  "

# Rows 48/49: RenderToken and BabyDogeCoin swap ranking positions, each with a freshly scraped volume figure
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.50"
$ws.Range("E48").Value = "  -4.40%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  +6.33%  "
